$wb = $excel.ActiveWorkbook

# Overview sheet: G2 (Latest HO Xliff Generate Date)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 19:09:39"

# zh-cn sheet: H2 (Correspond Handoff Datetime), K2 (Correspond Handback DateTime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 19:09:34"
$wsZhCn.Range("K2").Value = "2016-08-22 19:10:00"

# de-de sheet: H2 (Correspond Handoff Datetime), K2 (Correspond Handback DateTime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 19:09:39"
$wsDeDe.Range("K2").Value = "2016-08-22 19:10:23"
